$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing "Altro" (currently row 9) down to row 10
$ws.Rows.Item(9).Insert()

# New row 9: ID 8, "Residenza temporanea"
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Residenza temporanea"

# Row 10 now holds the old row-9 data (ID 9, "Altro") - already shifted by Insert()

# New row 11: ID 10, "Revisione onomastica stradale"
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Revisione onomastica stradale"

# Update selection to match target (activeCell D15, sqref D15)
$ws.Range("D15").Select()
